# The workbook's data rows (2..N) get re-sorted in ascending order by
# column A ("Id"). Row 1 is the header and stays in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data region starting below the header.
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
$lastCol = $ws.Cells(1, $ws.Columns.Count).End(-4159).Column

$dataRange = $ws.Range($ws.Cells(2, 1), $ws.Cells($lastRow, $lastCol))
$sortKey = $ws.Range($ws.Cells(2, 1), $ws.Cells($lastRow, 1))

# Sort ascending (1 = xlAscending) by the Id column, no header in range.
$dataRange.Sort($sortKey, 1)
